$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.929.88'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.811.63'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.52'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6080'
$ws.Range("E6").Value = '  -2.95%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9995'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07275'
$ws.Range("E8").Value = '  -2.89%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2861'
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.77'
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07631'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '1.798.37'
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.918'
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6573'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '80.87'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008901'
$ws.Range("E16").Value = '  -4.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.836'
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("D18").Value = '28.880.37'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '2.036.95'
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.03'
$ws.Range("E20").Value = '  +6.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.37'
$ws.Range("E21").Value = '  -1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9990'
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.077'
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9994'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.00'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1398'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.373'
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.54'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05563'
$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.068'
$ws.Range("E31").Value = '  +0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.065'
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.810'
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7305'
$ws.Range("E35").Value = '  -2.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.126'
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.619'
$ws.Range("E37").Value = '  -1.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.802'
$ws.Range("E38").Value = '  +1.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01746'
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("D40").Value = '1.190.17'
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.332'
$ws.Range("E41").Value = '  -2.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8886'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9991'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.43'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").Value = '1.937.48'
$ws.Range("E45").Value = '  -2.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000122'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '63.99'
$ws.Range("E47").Value = '  -1.94%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5063'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3961'
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.951'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05767'
$ws.Range("E51").Value = '  -0.91%  '
